$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$preValue = 0.625
$recallValue = 0.08196721311475409
$f1Value = 0.144927536231884

$startRow = 2
$endRow = 100

for ($row = $startRow; $row -le $endRow; $row++) {
    $ws.Cells.Item($row, 2).Value = $preValue
    $ws.Cells.Item($row, 3).Value = $recallValue
    $ws.Cells.Item($row, 4).Value = $f1Value
}
